$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A date cells to remain plain text so Excel doesn't
# auto-convert the dash-separated strings into date serial numbers.
$ws.Range("A3:A21").NumberFormat = "@"

# Update date strings (slash -> dash separators), keep as text
$ws.Range("A3").Value  = "28-07-2022"
$ws.Range("A4").Value  = "01-08-2022"
$ws.Range("A5").Value  = "04-08-2022"
$ws.Range("A6").Value  = "08-08-2022"
$ws.Range("A7").Value  = "11-08-2022"
$ws.Range("A8").Value  = "15-08-2022"
$ws.Range("A9").Value  = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# Update attendance counts that changed
$ws.Range("D3").Value  = 1
$ws.Range("G3").Value  = 1

$ws.Range("D4").Value  = 1
$ws.Range("G4").Value  = 1

$ws.Range("D6").Value  = 1
$ws.Range("E6").Value  = 1
$ws.Range("H6").Value  = 0

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0
